# Rename/shorten a handful of BART station names on the "Station Name"
# column (column B) of the station-codes sheet. Each row is located by its
# two-letter (or numeric) station code in column A so the edit is robust to
# the exact row position.
#
#   EN  El Cerrito Del Norte                -> El Cerrito del Norte
#   BK  Berkeley                            -> Downtown Berkeley
#   19  19th Street Oakland                 -> 19th Street
#   12  12th Street / Oakland City Center   -> 12th Street
#   BF  Bayfair                             -> Bay Fair
#   ED  Dublin/Pleasanton                   -> Dublin
#   WP  Pittsburg/Bay Point                 -> Pittsburg
#   SO  San Francisco International Airport -> SFO
#   WD  West Dublin/Pleasanton              -> West Dublin
#   OA  Oakland International Airport       -> OAK
#   BE  Berryessa / North San José          -> Berryessa

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$codes = $ws.Columns("A")

$renames = @{
    "EN" = "El Cerrito del Norte";
    "BK" = "Downtown Berkeley";
    "19" = "19th Street";
    "12" = "12th Street";
    "BF" = "Bay Fair";
    "ED" = "Dublin";
    "WP" = "Pittsburg";
    "SO" = "SFO";
    "WD" = "West Dublin";
    "OA" = "OAK";
    "BE" = "Berryessa";
}

# Fallback row numbers (1-based) in case a code can't be located via Find,
# matching the original layout of the sheet.
$fallbackRows = @{
    "EN" = 3;
    "BK" = 6;
    "19" = 9;
    "12" = 10;
    "BF" = 15;
    "ED" = 38;
    "WP" = 40;
    "SO" = 43;
    "WD" = 45;
    "OA" = 46;
    "BE" = 51;
}

foreach ($code in $renames.Keys) {
    $newName = $renames[$code]
    $cell = $codes.Find($code)
    if ($cell -ne $null) {
        $row = $cell.Row
    } else {
        $row = $fallbackRows[$code]
    }
    $ws.Cells.Item($row, 2).Value = $newName
}
